# Auto-generated: scheduled-runner price refresh for Seraph_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H..N) per sheet
# with freshly pulled market data. One cell (CRP!N141) had no HQ leve price
# to compare against, so its profit cell is cleared rather than written.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 392.875
$ws.Range("I9").Value = 267.5
$ws.Range("J9").Value = 518.25
$ws.Range("K9").Value = 267.5
$ws.Range("L9").Value = 518.25
$ws.Range("M9").Value = -98.5
$ws.Range("N9").Value = -856.25
$ws.Range("H15").Value = 667.35297
$ws.Range("I15").Value = 667.35297
$ws.Range("K15").Value = 2002.05891
$ws.Range("M15").Value = -1833.05891
$ws.Range("H33").Value = 164.6923
$ws.Range("I33").Value = 164.6923
$ws.Range("K33").Value = 164.6923
$ws.Range("M33").Value = 64.30770000000001
$ws.Range("H43").Value = 8499
$ws.Range("J43").Value = 8499
$ws.Range("L43").Value = 8499
$ws.Range("N43").Value = -8637
$ws.Range("H92").Value = 380.14285
$ws.Range("I92").Value = 380.14285
$ws.Range("K92").Value = 380.14285
$ws.Range("M92").Value = 867.85715
$ws.Range("H138").Value = 3965.476
$ws.Range("I138").Value = 1265.3334
$ws.Range("K138").Value = 3796.0002
$ws.Range("M138").Value = 1343.9998
$ws.Range("H141").Value = 6111
$ws.Range("I141").Value = 5314.6665
$ws.Range("K141").Value = 15943.9995
$ws.Range("M141").Value = -10763.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4166.6665
$ws.Range("I2").Value = 3750
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 3750
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -3637
$ws.Range("N2").Value = -5226
$ws.Range("H5").Value = 982.5
$ws.Range("I5").Value = 980
$ws.Range("K5").Value = 980
$ws.Range("M5").Value = -868
$ws.Range("H32").Value = 11332.179
$ws.Range("I32").Value = 7605.1816
$ws.Range("K32").Value = 7605.1816
$ws.Range("M32").Value = -7318.1816
$ws.Range("H63").Value = 5443
$ws.Range("J63").Value = 8501
$ws.Range("L63").Value = 8501
$ws.Range("N63").Value = -9873
$ws.Range("H66").Value = 5443
$ws.Range("J66").Value = 8501
$ws.Range("L66").Value = 42505
$ws.Range("N66").Value = -49369
$ws.Range("H116").Value = 4166.6665
$ws.Range("I116").Value = 3750
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 3750
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1456
$ws.Range("N116").Value = -9588
$ws.Range("H122").Value = 502905.2
$ws.Range("I122").Value = 1112622.6
$ws.Range("K122").Value = 3337867.8
$ws.Range("M122").Value = -3335417.8
$ws.Range("H132").Value = 1320.75
$ws.Range("I132").Value = 1421
$ws.Range("J132").Value = 1020
$ws.Range("K132").Value = 4263
$ws.Range("L132").Value = 3060
$ws.Range("M132").Value = -1733
$ws.Range("N132").Value = -8120

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4166.6665
$ws.Range("I3").Value = 3750
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 3750
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -3636
$ws.Range("N3").Value = -5228
$ws.Range("H4").Value = 982.5
$ws.Range("I4").Value = 980
$ws.Range("K4").Value = 980
$ws.Range("M4").Value = -865

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 256.625
$ws.Range("I10").Value = 283.83334
$ws.Range("K10").Value = 283.83334
$ws.Range("M10").Value = -144.83334
$ws.Range("H31").Value = 5340.8184
$ws.Range("I31").Value = 2098.3333
$ws.Range("J31").Value = 6556.75
$ws.Range("K31").Value = 2098.3333
$ws.Range("L31").Value = 6556.75
$ws.Range("M31").Value = -1803.3333
$ws.Range("N31").Value = -7146.75
$ws.Range("H34").Value = 5340.8184
$ws.Range("I34").Value = 2098.3333
$ws.Range("J34").Value = 6556.75
$ws.Range("K34").Value = 2098.3333
$ws.Range("L34").Value = 6556.75
$ws.Range("M34").Value = -1896.3333
$ws.Range("N34").Value = -6960.75
$ws.Range("H58").Value = 3798.4092
$ws.Range("I58").Value = 1308.7778
$ws.Range("K58").Value = 1308.7778
$ws.Range("M58").Value = -1105.7778
$ws.Range("H99").Value = 15593.044
$ws.Range("I99").Value = 13248.889
$ws.Range("J99").Value = 17100
$ws.Range("K99").Value = 13248.889
$ws.Range("L99").Value = 17100
$ws.Range("M99").Value = -11750.889
$ws.Range("N99").Value = -20096
$ws.Range("H122").Value = 2246.7727
$ws.Range("I122").Value = 2318.95
$ws.Range("K122").Value = 6956.849999999999
$ws.Range("M122").Value = -4506.849999999999
$ws.Range("H126").Value = 15593.044
$ws.Range("I126").Value = 13248.889
$ws.Range("J126").Value = 17100
$ws.Range("K126").Value = 39746.667
$ws.Range("L126").Value = 51300
$ws.Range("M126").Value = -37276.667
$ws.Range("N126").Value = -56240
$ws.Range("H132").Value = 3610
$ws.Range("I132").Value = 3610
$ws.Range("K132").Value = 10830
$ws.Range("M132").Value = -8300
$ws.Range("H134").Value = 3262.6843
$ws.Range("I134").Value = 2417.75
$ws.Range("J134").Value = 4711.143
$ws.Range("K134").Value = 7253.25
$ws.Range("L134").Value = 14133.429
$ws.Range("M134").Value = -4718.25
$ws.Range("N134").Value = -19203.429
$ws.Range("H136").Value = 3798.4092
$ws.Range("I136").Value = 1308.7778
$ws.Range("K136").Value = 3926.3334
$ws.Range("M136").Value = -1376.3334
$ws.Range("H141").Value = 105275
$ws.Range("I141").Value = 105275
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 105275
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -100095
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 55607.11
$ws.Range("I2").Value = 83360.5
$ws.Range("J2").Value = 100.333336
$ws.Range("K2").Value = 500163
$ws.Range("L2").Value = 602.000016
$ws.Range("M2").Value = -500050
$ws.Range("N2").Value = -828.000016
$ws.Range("H11").Value = 785.7143
$ws.Range("I11").Value = 250
$ws.Range("K11").Value = 750
$ws.Range("M11").Value = -610

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6935.1113
$ws.Range("I70").Value = 5807
$ws.Range("K70").Value = 5807
$ws.Range("M70").Value = -5537
$ws.Range("H73").Value = 6935.1113
$ws.Range("I73").Value = 5807
$ws.Range("K73").Value = 5807
$ws.Range("M73").Value = -4871
$ws.Range("H122").Value = 65295.75
$ws.Range("I122").Value = 2548.6924
$ws.Range("K122").Value = 7646.0772
$ws.Range("M122").Value = -5196.0772
$ws.Range("H132").Value = 4778.2
$ws.Range("I132").Value = 2751.3333
$ws.Range("K132").Value = 8253.999899999999
$ws.Range("M132").Value = -5723.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1933
$ws.Range("I40").Value = 1499.5
$ws.Range("K40").Value = 1499.5
$ws.Range("M40").Value = -1363.5
$ws.Range("H61").Value = 13999
$ws.Range("I61").Value = 13999
$ws.Range("K61").Value = 13999
$ws.Range("M61").Value = -13797
$ws.Range("H113").Value = 13999
$ws.Range("I113").Value = 13999
$ws.Range("K113").Value = 13999
$ws.Range("M113").Value = -11829
$ws.Range("H122").Value = 7397.1177
$ws.Range("I122").Value = 7988.727
$ws.Range("J122").Value = 6312.5
$ws.Range("K122").Value = 23966.181
$ws.Range("L122").Value = 18937.5
$ws.Range("M122").Value = -21516.181
$ws.Range("N122").Value = -23837.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5802.2
$ws.Range("I62").Value = 3928.923
$ws.Range("J62").Value = 7234.706
$ws.Range("K62").Value = 3928.923
$ws.Range("L62").Value = 7234.706
$ws.Range("M62").Value = -3304.923
$ws.Range("N62").Value = -8482.706
$ws.Range("H65").Value = 5802.2
$ws.Range("I65").Value = 3928.923
$ws.Range("J65").Value = 7234.706
$ws.Range("K65").Value = 19644.615
$ws.Range("L65").Value = 36173.53
$ws.Range("M65").Value = -16524.615
$ws.Range("N65").Value = -42413.53
$ws.Range("H107").Value = 429.10526
$ws.Range("I107").Value = 342.6
$ws.Range("J107").Value = 460
$ws.Range("K107").Value = 1027.8
$ws.Range("L107").Value = 1380
$ws.Range("M107").Value = 892.1999999999998
$ws.Range("N107").Value = -5220
$ws.Range("H136").Value = 1532.3096
$ws.Range("I136").Value = 1167.6316
$ws.Range("K136").Value = 3502.8948
$ws.Range("M136").Value = -952.8948
